$wb = $excel.ActiveWorkbook

$wsDaily = $wb.Worksheets.Item("Daily")
$wsHourly = $wb.Worksheets.Item("Hourly")

# --- Daily sheet, row 2 ---
$wsDaily.Range("G2").Value = 3190.91
$wsDaily.Range("H2").Value = 6465.48
$wsDaily.Range("I2").Value = 765.74
$wsDaily.Range("J2").Value = 3190.91
$wsDaily.Range("K2").Value = 6145.8
$wsDaily.Range("L2").Value = 768.09

# --- Hourly sheet ---
# Row 9
$wsHourly.Range("L9").Value = 16.07

# Row 10
$wsHourly.Range("I10").Value = 445.94
$wsHourly.Range("L10").Value = 409.19

# Row 12
$wsHourly.Range("I12").Value = 744.02

# Row 13
$wsHourly.Range("I13").Value = 792.87
$wsHourly.Range("L13").Value = 768.64

# Row 14
$wsHourly.Range("H14").Value = 513.35
$wsHourly.Range("I14").Value = 810.58
$wsHourly.Range("K14").Value = 513.35
$wsHourly.Range("L14").Value = 781.95
$wsHourly.Range("M14").Value = 106.44

# Row 15
$wsHourly.Range("H15").Value = 495.15
$wsHourly.Range("I15").Value = 801.96

# Row 16
$wsHourly.Range("H16").Value = 424.38
$wsHourly.Range("I16").Value = 764.6799999999999
$wsHourly.Range("K16").Value = 424.38
$wsHourly.Range("L16").Value = 745.63

# Row 17
$wsHourly.Range("I17").Value = 687.15
$wsHourly.Range("J17").Value = 79.48
$wsHourly.Range("L17").Value = 669.01

# Row 18
$wsHourly.Range("H18").Value = 165.07
$wsHourly.Range("I18").Value = 533.08
$wsHourly.Range("K18").Value = 165.07
$wsHourly.Range("L18").Value = 509.24
$wsHourly.Range("M18").Value = 55.75

# Row 19
$wsHourly.Range("I19").Value = 177.07
$wsHourly.Range("L19").Value = 115
